$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column before column B. This shifts the existing
# "DevSuite" column from B to C and leaves a blank column B behind -
# column A ("ProdSuite") keeps its original formatting/width untouched.
$ws.Columns("B:B").Insert()

# Header row: new column B re-uses the old "ProdSuite" header that used
# to live in A1, and A1 becomes the new "Suite" label.
$oldA1 = $ws.Range("A1").Value2
$ws.Range("B1").Value2 = $oldA1
$ws.Range("A1").Value2 = "Suite"

# Data rows: column B is a duplicate of column A's environment name,
# column C keeps the values that used to be in column B (already shifted
# there automatically by the Insert above).
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 2).Value2 = $ws.Cells.Item($r, 1).Value2
}

# Match the new column B's width to column A's (as closely as the host
# allows) so the two "Suite" columns line up visually; column A and C
# retain their original widths since we never touch them directly.
$ws.Columns("B:B").ColumnWidth = 39
